$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report title (remove the period after "6.4.1.2") in Russian and English headers
$ws.Range("B1").Value = "6.4.1.2 Потери воды при транспортировке"
$ws.Range("C1").Value = "6.4.1.2 Percentage of water loss during transportation"

# Update the 2022 data values that changed
$ws.Range("P5").Value = 2388
$ws.Range("P10").Value = 335.3
$ws.Range("P16").Value = 27.3
$ws.Range("P21").Value = 24.3

# Update the active cell selection
$null = $ws.Range("S3").Select()
